# Scen_RES_SHARE_50%_24_7.xlsx — expand the single ANNUAL TFM_INS row into a
# full set of seasonal/day-night TimeSlice rows (WIN/SPR/SUM/AUT x DAY/NITE)
# repeated for each milestone year 2025..2050.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timeSlices = @("WIN_DAY", "WIN_NITE", "SPR_DAY", "SPR_NITE", "SUM_DAY", "SUM_NITE", "AUT_DAY", "AUT_NITE")
$years = @(2025, 2030, 2035, 2040, 2045, 2050)

# Common column values shared by every data row (B varies by timeslice, F by year).
$csetCn = "ELC_FIN_DEM"
$psetPn = "ELC_GRID_RES"
$attribute = "FLO_SHAR"
$limType = "FX"
$value = 0.5

$row = 6
foreach ($year in $years) {
    foreach ($ts in $timeSlices) {
        $ws.Cells.Item($row, 2).Value = $ts
        $ws.Cells.Item($row, 3).Value = $csetCn
        $ws.Cells.Item($row, 4).Value = $psetPn
        $ws.Cells.Item($row, 5).Value = $attribute
        $ws.Cells.Item($row, 6).Value = $year
        $ws.Cells.Item($row, 7).Value = $limType
        $ws.Cells.Item($row, 8).Value = $value
        $row = $row + 1
    }
}

# Reflect the final selection left behind in the saved workbook.
$ws.Range("B6:H53").Select()
